# Update "想去人数" (F column) figures for the 展览 and 全部类型 sheets.
# Both sheets contain the same event listing and therefore receive the
# identical set of updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1084
    "F3"  = 789
    "F5"  = 40
    "F8"  = 1924
    "F9"  = 7135
    "F10" = 896
    "F11" = 400
    "F12" = 331
    "F13" = 120
    "F14" = 390
    "F16" = 7066
    "F17" = 289
    "F18" = 1322
    "F19" = 145
    "F21" = 227
    "F22" = 132
    "F23" = 292
    "F24" = 128
    "F29" = 405
    "F30" = 605
    "F32" = 88
    "F33" = 51
    "F34" = 69
    "F36" = 71
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
